# MAJ de certaines images
# Adds numbered "ellipse" call-out bubbles (1/2/3) to slide 1 and
# (1/2) plus a new rounded-rectangle highlight box to slide 2, and
# tweaks the corner radius ("adj" handle) of the existing rounded
# rectangle on slide 2.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 1 (sldId 257) - add three numbered ellipse call-outs
# ---------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# The slide previously had a picture (id 2) deleted from it, so the
# next free shape id is 2, not 8. Real PowerPoint never re-issues a
# retired shape id, so "burn" the free id with a throw-away shape
# before adding the real content -- this reproduces the id 8/9/10
# numbering the authoring app produced.
$filler = $s1.Shapes.AddShape(9, 0, 0, 1, 1)
$filler.Delete()

# Use one of the existing rounded-rectangle shapes as a style donor:
# duplicating keeps the theme-linked <p:style> block, body/paragraph
# formatting, etc. that a brand-new AddShape() would not include.
$styleDonor1 = $s1.Shapes.Item(2)

$e1 = $styleDonor1.Duplicate().Item(1)
$e1.Name = "Ellipse 7"
$e1.AutoShapeType = 9
$e1.Left = 450.2450409700788
$e1.Top = 84.12803269606299
$e1.Width = 29.75496062992126
$e1.Height = 29.75496062992126
$e1.Line.Visible = $false
$e1.TextFrame.TextRange.Text = "1"
$e1.TextFrame.TextRange.LanguageID = "fr-CA"

$e2 = $styleDonor1.Duplicate().Item(1)
$e2.Name = "Ellipse 8"
$e2.AutoShapeType = 9
$e2.Left = 631.8133858267717
$e2.Top = 162.68134308267716
$e2.Width = 29.75496062992126
$e2.Height = 29.75496062992126
$e2.Line.Visible = $false
$e2.TextFrame.TextRange.Text = "2"
$e2.TextFrame.TextRange.LanguageID = "fr-CA"

$e3 = $styleDonor1.Duplicate().Item(1)
$e3.Name = "Ellipse 9"
$e3.AutoShapeType = 9
$e3.Left = 828.5043307086614
$e3.Top = 86.74622047244094
$e3.Width = 29.75496062992126
$e3.Height = 29.75496062992126
$e3.Line.Visible = $false
$e3.TextFrame.TextRange.Text = "3"
$e3.TextFrame.TextRange.LanguageID = "fr-CA"

# ---------------------------------------------------------------
# Slide 2 (sldId 258) - bump the rounded rectangle's corner radius,
# add a second rounded-rectangle highlight box and two numbered
# ellipse call-outs
# ---------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$existingRect = $s2.Shapes.Item(2)
$existingRect.Adjustments.Item(1) = 0.06423

$r1 = $existingRect.Duplicate().Item(1)
$r1.Name = "Rectangle : coins arrondis 3"
$r1.Left = 317.018031496063
$r1.Top = 133.76007874015747
$r1.Width = 151.63645669291338
$r1.Height = 17.75488188976378

$e4 = $existingRect.Duplicate().Item(1)
$e4.Name = "Ellipse 4"
$e4.AutoShapeType = 9
$e4.Left = 453.77700787401574
$e4.Top = 127.76
$e4.Width = 29.75496062992126
$e4.Height = 29.75496062992126
$e4.Line.Visible = $false
$e4.TextFrame.TextRange.Text = "1"
$e4.TextFrame.TextRange.LanguageID = "fr-CA"

$e5 = $existingRect.Duplicate().Item(1)
$e5.Name = "Ellipse 5"
$e5.AutoShapeType = 9
$e5.Left = 552.1770324740157
$e5.Top = 176.51055148110237
$e5.Width = 29.75496062992126
$e5.Height = 29.75496062992126
$e5.Line.Visible = $false
$e5.TextFrame.TextRange.Text = "2"
$e5.TextFrame.TextRange.LanguageID = "fr-CA"
